$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headline = 'miniso group holding limited unsponsored adr nysemnso given consensus recommendation of moderate buy by brokerage cleanspark nasdaqclsk research coverage started at ladenburg thalmsh sh i got accepted to harvard stanford and princeton thing my parent did at an early age it made all the difference tva report billion in operating revenue up from best free machine learning course on udemy youdao to report second quarter financial result on august zto to announce second quarter and halfyear interim financial result of on august u eastern time from the financial service csuite to leading a national school choice movement with erika donalds from the corner office closing the prime age employment gap netease to report second quarter financial result on august in remittix ethereumbased token might be your fast track to figure huge bonus are powering this surge cwm llc ha stake in cooperstandard holding inc nysecps invesco pharmaceutical etf nysearcapjp share sold by envestnet asset management inc ring energy inc nyseamericanrei share sold by bank of new york mellon corp intrinsic value partner llc ha million stock holding in ishares core sp smallcap etf nysearcaijr wealthquest corp sell share of ishares core sp smallcap etf nysearcaijr eastern bank sell share of ishares msci eafe etf nysearcaefa alexander randolph advisory inc buy share of ishares msci eafe etf nysearcaefa cordatus wealth management llc decrease stock position in ishares msci eafe etf nysearcaefa rvw wealth llc lower stock position in ishares msci eafe etf nysearcaefa u financial advisor llc grows stock holding in ishares msci eafe etf nysearcaefa tarbox family office inc sell share of ishares msci eafe etf nysearcaefa efg asset management america corp acquires share of ishares msci eafe etf nysearcaefa sound view wealth advisor group llc cut stock position in ishares russell growth etf nysearcaiwf ishares msci eafe etf nysearcaefa stock holding trimmed by lyell wealth management lp tillman hartley llc raise stock holding in ishares msci eafe etf nysearcaefa harvest fund management co ltd take position in eastgroup property inc nyseegp empirical financial service llc dba empirical wealth management ha million stock holding in ishares core sp smallcap etf nysearcaijr north star asset management inc acquires share of ishares msci eafe etf nysearcaefa efg asset management america corp sell share of ishares core sp smallcap etf nysearcaijr first national trust co ha million stock holding in ishares russell growth etf nysearcaiwf fiera capital corp ha million stake in tencent music entertainment group sponsored adr nysetme fiera capital corp sell share of blackrock muniholdings fund inc nysemhd envestnet asset management inc lower stock holding in invesco pharmaceutical etf nysearcapjp foundation investment advisor llc raise position in fidelity quality factor etf nysearcafqal foundation investment advisor llc decrease stake in global x sp catholic value etf nasdaqcath cwm llc buy share of cooperstandard holding inc nysecps fiera capital corp lower stake in inv vk mun opp nysevmo foundation investment advisor llc sell share of innovator u equity power buffer etf september batspsep foundation investment advisor llc purchase share of ishares russell midcap etf nysearcaiwr foundation investment advisor llc trim holding in ishares dow jones u etf nysearcaiyy foundation investment advisor llc buy share of innovator u equity ultra buffer etf september batsusep foundation investment advisor llc lower position in proshares ultra qqq nysearcaqld foundation investment advisor llc acquires share of american century u quality growth etf nysearcaqgro wealthquest corp ha million stock position in ishares core sp smallcap etf nysearcaijr empirical financial service llc dba empirical wealth management ha million stock holding in ishares core sp smallcap etf nysearcaijr alexander randolph advisory inc purchase share of ishares msci eafe etf nysearcaefa efg asset management america corp grows stake in ishares msci eafe etf nysearcaefa inmune bio inmb to release quarterly earnings on thursday rvw wealth llc sell share of ishares msci eafe etf nysearcaefa gen digital gen expected to announce earnings on thursday anghami nasdaqangh share set to reverse split on monday august th barclays forecast strong price appreciation for rtx nysertx stock stifel nicolaus issue pessimistic forecast for revvity nyservty stock price whitefiber inc plan million ipo for week of august th wyfi whitefiber inc plan million ipo for week of august th wyfi northrop grumman nysenoc given new price target at barclays apple doe it again promotes samsungs galaxy z flip by mistake sarepta therapeutic nasdaqsrpt raised to equal weight at barclays raymond james financial reaffirms outperform rating for alerus financial nasdaqalrs keycorp issue positive forecast for cadence design system nasdaqcdns stock price todd m schneider sell share of cintas corporation nasdaqctas stock columbia banking system nasdaqcolb given new price target at barclays jpmorgan chase co boost cadence design system nasdaqcdns price target to justin rosenstein sell share of asana inc nyseasan stock insider selling asana inc nyseasan director sell in stock east west bancorp inc nasdaqewbc ceo sell in stock holding ltd v sell share of cipher mining inc nasdaqcifr stock august j troendle sell share of medpace holding inc nasdaqmedp stock kellanova nysek major shareholder kellogg w k foundation trust sell share keefe bruyette wood cut meritage home nysemth price target to skywest nasdaqskyw stock price expected to rise evercore isi analyst say carvana co nysecvna ceo sell in stock insider selling medpace holding inc nasdaqmedp cfo sell share of stock michael o mccarthy iii sell share of quantumscape corporation nyseqs stock insider selling arista network inc nyseanet ceo sell share of stock august j troendle sell share of medpace holding inc nasdaqmedp stock keefe bruyette wood issue pessimistic forecast for meritage home nysemth stock price ubs group cut charter communication nasdaqchtr price target to john s jenkins sell share of te connectivity ltd nysetel stock tony xu sell share of doordash inc nasdaqdash stock vladimir galkin acquires share of newegg commerce inc nasdaqnegg stock tractor supply company nasdaqtsco cfo sell in stock doordash inc nasdaqdash ceo tony xu sell share tmobile u inc nasdaqtmus director telekom ag deutsche sell share tmobile u inc nasdaqtmus director telekom ag deutsche sell share peter l briger jr purchase share of microstrategy incorporated nasdaqmstr stock barclays cut hca healthcare nysehca price target to vanguard u momentum factor etf batsvfmo share purchased by empirical financial service llc dba empirical wealth management ishares russell etf nysearcaiwb share purchased by empirical financial service llc dba empirical wealth management synovus financial nysesnv price target cut to by analyst at stephen stephen raise ryder system nyser price target to ishares msci usa quality factor etf batsqual is empirical financial service llc dba empirical wealth management th largest position scotiabank increase lightspeed commerce nyselspd price target to range resource nyserrc price target raised to alexander randolph advisory inc sell share of royce otc micro nysermt ubs group raise u food nyseusfd price target to empirical financial service llc dba empirical wealth management grows position in vanguard smallcap value etf nysearcavbr empirical financial service llc dba empirical wealth management grows stock position in vanguard ftse allworld exus smallcap etf nysearcavss barrington research reaffirms outperform rating for quad graphic nysequad bg food inc nysebgs given average recommendation of reduce by analyst the real genocide isnt happening in gaza wyoming parent shouldnt lose hope about paused education choice fund republican already seem to be surrendering north carolina u senate race opinion alvia lewis frey despite friend effort i still forget thing what would carson say buuu group limited buuu plan to raise million in august th ipo barclays raise rtx nysertx price target to corporacion inmobiliaria vesta nysevtmx stock price expected to rise barclays analyst say trump warns canada backing palestinian state could damage u trade deal ptc therapeutic nasdaqptct price target raised to holding ltd v sell share of cipher mining inc nasdaqcifr stock teledyne technology incorporated nysetdy director sell in stock dominic ng sell share of east west bancorp inc nasdaqewbc stock insider selling maplebear inc nasdaqcart insider sell share of stock michael o mccarthy iii sell share of quantumscape corporation nyseqs stock cadence design system nasdaqcdns given new price target at jpmorgan chase co hca healthcare nysehca given a price target at royal bank of canada insider selling medpace holding inc nasdaqmedp president sell share of stock coliseum capital management l acquires share of sonos inc nasdaqsono stock ishares esg aware msci usa etf nasdaqesgu holding cut by empirical financial service llc dba empirical wealth management hershey nysehsy stock price expected to rise stifel nicolaus analyst say vanguard smallcap value etf nysearcavbr share acquired by empirical financial service llc dba empirical wealth management ishares short maturity municipal bond active etf batsmear share acquired by alexander randolph advisory inc standex international nysesxi receives outperform rating from barrington research ubs group issue positive forecast for hexcel nysehxl stock price u food nyseusfd price target raised to citigroup ha lowered expectation for waste connection nysewcn stock price procore technology nysepcor price target raised to lithuania prime minister gintautas paluckas step down after investigation and protest calling for his resignation lithuania prime minister gintautas paluckas step down after investigation and protest calling for his resignation lithuania prime minister gintautas paluckas step down after investigation and protest calling for his resignation republican in texas release a proposal for a new congressional map brett j hart sell share of united airline holding inc nasdaqual stock trump medium technology group corp nasdaqdjtww short interest down in july ffbw inc nasdaqffbw see large decline in short interest kamala harris political future short interest in great elm capital corp note due nasdaqgeccz decline by aureus greenway holding inc nasdaqagh major shareholder sell in stock fidelity disruptive finance etf nasdaqfdff trading down time to sell atlas nasdaqatcol share up should you buy global gas co nasdaqhgasw short interest update liontown resource otcmktslinrf stock price down here what happened virnetx nysevhc share cross above day moving average here why senator alex padilla on redistricting effort by republican sen alex padilla say republican want to rig midterm with redistricting plan global x interest rate hedge etf nysearcarate stock price up here what happened cardio diagnostics nasdaqcdio trading higher time to buy benev capital bevff to release earnings on thursday avant brand otcmktsavtbf stock price down here why wearable device nasdaqwlds trading down time to sell coinshares international otcmktscnsrf trading higher still a buy first mid bancshares fmbh to release quarterly earnings on thursday heritage global hgbl to release earnings on thursday ecovyst ecvt projected to post earnings on thursday local leader to face texas lawmaker today a question linger about deadly july flooding trump blast failed former russian president medvedev to watch his word over war threat methuen mayor on publicly supporting market basket ceo if it not broken why fix it the democratic party face an uphill battle ahead of the next election cycle the democratic party is still trying to figure out it message after s defeat hoxton planning management llc acquires share of ishares global infrastructure etf nasdaqigf schwab shortterm u treasury etf nysearcascho share sold by towercrest capital management atlas wealth llc ha million stake in john hancock multifactor mid cap etf nysearcajhmm grantvest financial group llc raise holding in vanguard midcap value etf nysearcavoe grantvest financial group llc raise position in vanguard russell etf nasdaqvtwo grantvest financial group llc ha million stock holding in first trust rising dividend achiever etf nasdaqrdvy atlas wealth llc decrease holding in invesco sp smallcap financials etf nasdaqpscf cordatus wealth management llc raise position in hsbc holding plc nysehsbc lyell wealth management lp ha million position in ishares ibonds dec term corporate etf nysearcaibdu tillman hartley llc sell share of vanguard shortterm treasury etf nasdaqvgsh no jail time for man who vowed to kill a many jew a possible in toronto barclays forecast strong price appreciation for cemex nysecx stock tigress financial raise carnival nyseccl price target to d boral capital reiterates buy rating for canfite biopharma nysecanf north west tsenwc stock rating upgraded by royal bank of canada first quantum mineral tsefm price target raised to c d boral capital reaffirms buy rating for eon resource nyseamericaneonr keefe bruyette wood increase seacoast banking corporation of florida nasdaqsbcf price target to stephen issue positive forecast for autonation nysean stock price labor union representing some costco employee for short crossword clue im in no mood to argue with you crossword clue collaborative wealth managment inc decrease stock holding in spdr portfolio short term treasury etf nysearcaspts collaborative wealth managment inc sell share of spdr blackstone senior loan etf nysearcasrln lyell wealth management lp raise holding in vaneck merk gold etf nysearcaounz invesco national amtfree municipal bond etf nysearcapza stock position raised by collaborative wealth managment inc collaborative wealth managment inc ha position in proshares sp aristocrat etf batsnobl goodman financial corp grows stake in spdr portfolio sp mid cap etf nysearcaspmd collaborative wealth managment inc ha stock position in first trust tcw opportunistic fixed income etf nasdaqfixd collaborative wealth managment inc cut position in ishares core dividend growth etf nysearcadgro collaborative wealth managment inc buy share of spdr portfolio developed world exus etf nysearcaspdw vanguard shortterm corporate bond etf nasdaqvcsh share sold by goodman financial corp microsoft and openai try to rebalance their relationship redburn atlantic upgrade moonlake immunotherapeutics nasdaqmltx to hold neuren pharmaceutical otcmktsnurpf trading up here what happened gentex nasdaqgntx price target raised to singing machine nasdaqmics trading down here why domino pizza group otcmktsdmpzf share down should you sell agronomics otcmktsagnmf trading down whats next life nasdaqlif stock price expected to rise stifel nicolaus analyst say fortinet nasdaqftnt given new price target at well fargo company global oncology diagnostics market to reach billion by u childhood vaccination rate fall again a exemption set another record u childhood vaccination rate fall again a exemption set another record u childhood vaccination rate fall again a exemption set another record scotland give goahead for world largest offshore wind farm in israel there is a sea change a former official see opposition growing u kindergarten vaccination rate inched down again last year and the exemption rate rose again u kindergarten vaccination rate inched down again last year and the exemption rate rose again state street new coveredcall etf could turbocharge your monthly income what copper collapse mean for gold silver and mining stock gotb take over new york city with immersive curlpowered popups carronade share perspective on viasat tiktok mom hack tariff with backtoschool budget no idea if this is possible china is researching way to sabotage elon musk starlink satellite with laser sec launch project crypto here what it mean for you patent fee reimagined evaluating the trump administration valuedriven fee structure and it impact on innovator stock market news today market up on meta microsoft boost sp seeking alpha popular energy drink pulled from shelf for containing vodka collaborative wealth managment inc reduces stock position in the walt disney company nysedis bottling plant proposed in perth under planning board review stay adventure ready this summer with the igloo qt cooler now off employerfriendly noncompete legislation becomes law in florida erayak power solution group announces million registered direct offering lyno ai best presale to invest in bitcoin consolidates after fomc meeting what doe powell decision mean for eth xrp dogecoin best cryptos to buy and hold top pick for and beyond campines agile navigation through the trade war effect gray medium name eric walter a general manager of knoe and kaqy in monroe louisiana best workout apps in rudolph libbe group top crains annual ne ohio contractor list eu eye group lng buy to meet u commitment nvrh cost containment step include reduction in staff position walmart is selling a sophisticated burberry cologne for just and shopper like it nice subtle scent the linkedin hiring playbook for manager adam emrich named vice president of customer solution and sale at southern linc highyield dividend giant demolished analyst q earnings expectation construction on million white house ballroom to begin in september construction on million white house ballroom to begin in september data scientist ibm a scam or an inside job how coindcx lost m in crypto survey reveals barrier in u cleantech innovation estonia emerges a strategic partner survey reveals barrier in u cleantech innovation estonia emerges a strategic partner tenable holding analyst raise their forecast after betterthanexpected q earnings letter dont use tax dollar to betray due process datarobot zillow la vega is short k home a housing deficit continues to grow where venetian elegance meet miami skyline the european allure of cipriani residence parkview health commits to building fullservice hospital in west lafayette higher sugar output in brazil weighs on price real estate expert lance kammes highlight why family are moving to wheaton il in hellonation magazine homebuying expert sue lantz explains mortgage preapprovals in hellonation magazine fence expert corin clark explains what people misunderstand about residential installation in hellonation magazine almost through phase thanks to coinmarketcap listing ruvi ai ruvi holder eye gain catch it before it outpaces cardano ada cocoa price settle sharply higher on supply woe expert say ruvi ai ruvi could be this year biggest gainer presale skyrocket after coinmarketcap listing just like early ethereum eth forget prime walmart ha jlab earbuds for under brazil tariff risk underpin arabica coffee price justice kavanaugh defends supreme court terse emergency docket order my turn august is national breastfeeding month forward together lowcost highyield invesco etf with strong track record power your adventure with off the samsung t shield tb portable ssd ethereum eye new rally nomy research identifies ai boom a key driver gold gann cycle suggests bottom is inbut resistance loom these analyst boost their forecast on trane technology following q result mark carney fast fact how new u tariff on the eu and uk could affect the cost of booze my turn the connecticut river eight project figma triple in firstday trading figmas stock soar in it highly anticipated ipo market cap instantly hit b techcrunch trump announces day negotiating period with mexico a tariff rate stay in place invent inventor develops central vacuum system for vehicle interior pbt what to say and not say during a job interview figma stock blast off with wellreceived ipo the latest judge question trump authority to impose sweeping tariff aug deadline loom the latest judge question trump authority to impose sweeping tariff aug deadline loom the latest judge question trump authority to impose sweeping tariff aug deadline loom the latest judge question trump authority to impose sweeping tariff aug deadline loom the latest judge question trump authority to impose sweeping tariff aug deadline loom apple q earnings today analyst want update on ai progress iphone demand with the stock down business insider kolter urban selects firstservice residential to manage selene oceanfront residence kolter urban selects firstservice residential to manage selene oceanfront residence ohfa expands team with diverse hiring spree to bolster housing effort in ohio figma more than triple in nyse debut after selling share at rail deal will open new market for top u container port yahoo finance rail deal will open new market for top u container port uk want lender to simplify digital loan application amazon ad business jump year over year making b spire declares dividend spire declares dividend a most california city fall fresno and bakersfield climb in renter ranking justin timberlake reveals lyme disease diagnosis flo reach settlement in class action lawsuit after judge said he would toss plantiffs claim time to pay up on time uk government announces crackdown on late payment to small business starbucks could see a breakdown after stalling around key level according to the chart global trade in photo a a new deadline approach on trump tariff global trade in photo a a new deadline approach on trump tariff global trade in photo a a new deadline approach on trump tariff i cant watch this candace parker call out wnba for demanding big pay despite terrible play ap business summarybrief at pm edt china summons nvidia over backdoor safety risk in h chip construction of trump million white house ballroom to begin in september live now schumer senate democrat discus the effect of trump tariff trump sends witkoff and huckabee to inspect food distribution in gaza in israel there is a sea change a former official see opposition growing trump sends witkoff and huckabee to inspect food distribution in gaza more than national guard troop leaving la newsom say trump political theater backfired addison primary election result bidens postpresidential fall from grace even sadder than you think st louis to restart water shutoffs ending covid pause u and ecuador sign agreement to combat transnational crime epic game store is bringing fortnite back to google play detroit mayor primary election result feminist rape hoax obliterates life of canadian hockey star a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration man indicted after string of burglary at business storage unit in multnomah village trump sends witkoff and huckabee to inspect food distribution in gaza judge overseeing trump case speak out against threat judge overseeing trump case speak out against threat judge overseeing trump case speak out against threat should yearolds have a vote editorial branighan branighan branighan branighan branighan branighan branighan branighan branighan branighan it field trip day at the retirement home lefty engage in geriatric theatrics at wi rep steils home judge overseeing trump case speak out against threat judge overseeing trump case speak out against threat illinois lawmaker official seek input into how auto insurance rate are set judge overseeing trump case speak out against threat trump go off on secondtier senator josh hawley for supporting bill he hate national view sorry big bird but you were the dumbest waste of tax dollar white house announces u official will travel into gaza lehigh county controller mark pinsley enters crowded race for congress in lehigh valley elephant in the room the latest america school child face trump order to resume pushup and pullups trump envoy to visit gaza a pressure mount on israel karoline leavitt announces white house will begin construction of new state ballroom ap news summary at pm edt the order of maine november ballot question ha been set justice kavanaugh defends supreme court terse emergency docket order justice kavanaugh defends supreme court terse emergency docket order fbi never probed confidential talk between dnc clinton campaign and george soros group bombshell intel file show mark carney fast fact froma harrop news for the downwardly mobile elite youre still elite epping police investigating biting incident at new england dragway a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration a look at college with federal money targeted by the trump administration the latest judge question trump authority to impose sweeping tariff aug deadline loom the latest judge question trump authority to impose sweeping tariff aug deadline loom the latest judge question trump authority to impose sweeping tariff aug deadline loom john brennan hope career dojfbi official are not going to carry on this absurd effort nudge nudge important to report traffic incident accident letter how trump made one of the most controversial post policy the new normal guest columnist question for rep comer klamath county board of commissioner meeting zoombombed what maine voter need to know about upcoming referendum question what maine voter need to know about upcoming referendum question justin timberlake reveals lyme disease diagnosis trump to approve new gaza aid plan after witkoff visit friday white house kamala harris decision kickstarts the democratic presidential primary takeaway how aid is distributed in gaza here the big beautiful bill million pet project that few noticed zuckerberg overtakes bezos for world thirdrichest person a meta rally add billion to his fortune the latest america school child face trump order to resume pushup and pullups the latest america school child face trump order to resume pushup and pullups trump special envoy witkoff and u ambassador huckabee will inspect food distribution in gaza trump special envoy witkoff and u ambassador huckabee will inspect food distribution in gaza trump special envoy witkoff and u ambassador huckabee will inspect food distribution in gaza inside trump huge new vanity project a he build ft gold ballroom exclusive ron johnson urge hegseth to address military covid mandate victim following dcnf report u special envoy ambassador will inspect food distribution in gaza donald trump joke on south park help lift viewership by white house call sydney sweeney ad outrage cancel culture run amok rnc chair michael whatley announces senate bid in north carolina justin early to run for house district university of iowa need a reminder that no one is above antidei law watch grant help develop housing in kensington project how the left lost the working class trump threatens putin ally after he dared to compare him to sleepy joe trump sending witkoff huckabee to gaza amid aid and hunger crisis'

$ws.Cells.Item(28, 1).Value = 45870
$ws.Cells.Item(28, 1).NumberFormat = $ws.Cells.Item(27, 1).NumberFormat
$ws.Cells.Item(28, 2).Value = 6238.009765625
$ws.Cells.Item(28, 3).Value = 6287.27978515625
$ws.Cells.Item(28, 4).Value = 6212.68994140625
$ws.Cells.Item(28, 5).Value = 6287.27978515625
$ws.Cells.Item(28, 6).Value = 5827150000
$ws.Cells.Item(28, 7).Value = -0.0159921331401484
$ws.Cells.Item(28, 8).Value = -1
$ws.Cells.Item(28, 9).Value = $headline
$ws.Cells.Item(28, 10).Value = $headline
$ws.Cells.Item(28, 11).Value = 0.096
$ws.Cells.Item(28, 12).Value = 0.05
$ws.Cells.Item(28, 13).Value = 0.854
$ws.Cells.Item(28, 14).Value = 0.9991
$ws.Cells.Item(28, 15).Value = 75
$ws.Cells.Item(28, 16).Value = 87
$ws.Cells.Item(28, 17).Value = 0
$ws.Cells.Item(28, 18).Value = 309
$ws.Cells.Item(28, 19).Value = 108
$ws.Cells.Item(28, 20).Value = 410
$ws.Cells.Item(28, 21).Value = 160
$ws.Cells.Item(28, 22).Value = 68
$ws.Cells.Item(28, 23).Value = 40
$ws.Cells.Item(28, 24).Value = 168

Write-Output "done"
